# Schedule all student successfully
# Rotates the StudentID values in column B (rows 2-38) - the first 14
# student IDs (rows 2-15) are moved to the end of the block (rows 25-38),
# and the remaining 23 (originally rows 16-38) shift up to rows 2-24.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$startRow = 2
$endRow = 38

# Capture current StudentID values (column B) for rows 2..38
$values = @()
for ($r = $startRow; $r -le $endRow; $r++) {
    $values += $ws.Cells.Item($r, 2).Value()
}

$count = $values.Count
$shift = 14  # number of rows to rotate by (rows 2-15 move to the tail)

$rotated = @()
for ($i = 0; $i -lt $count; $i++) {
    $srcIndex = ($i + $shift) % $count
    $rotated += $values[$srcIndex]
}

for ($i = 0; $i -lt $count; $i++) {
    $ws.Cells.Item($startRow + $i, 2).Value = $rotated[$i]
}

# Update the view state to match: scrolled a bit further down, new selection.
$ws.Application.ActiveWindow.ScrollRow = 27
$ws.Range("C36").Select()
